# Weekly fruit/vegetable data update: insert a new record as row 131
# (pushing the existing rows 131-147 down to 132-148) on the single
# worksheet of this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before the current row 131, shifting rows
# 131-147 down to 132-148 (row 147's data ends up as the new row 148).
$ws.Rows.Item(131).Insert()

# Populate the newly-inserted row 131 with the new weekly record.
$ws.Cells.Item(131, 1).Value = 8
$ws.Cells.Item(131, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(131, 3).Value = "Coquimbo"
$ws.Cells.Item(131, 4).Value = 44449
$ws.Cells.Item(131, 5).Value = 4
$ws.Cells.Item(131, 6).Value = 100112012
$ws.Cells.Item(131, 7).Value = "Espinaca"
$ws.Cells.Item(131, 8).Value = "Sin especificar"
$ws.Cells.Item(131, 9).Value = "Primera"
$ws.Cells.Item(131, 10).Value = 3000
$ws.Cells.Item(131, 11).Value = 400
$ws.Cells.Item(131, 12).Value = 500
$ws.Cells.Item(131, 13).Value = 450
$ws.Cells.Item(131, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(131, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(131, 16).Value = 900
$ws.Cells.Item(131, 17).Value = 0.5
$ws.Cells.Item(131, 18).Value = "Hortaliza"
